$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows right below the last existing data row (105). Inserting
# (rather than just writing past the end of the used range) makes Excel
# copy row 105's formatting down into the new rows, matching the carried
# over cell styles on columns A-H seen in the target edit.
$ws.Rows("106:111").Insert()

# --- Row 106 : 16-04-2018 ---
$ws.Range("A106").Value = "16-04-2018"
$ws.Range("B106").Value = 1994
$ws.Range("C106").Value = 58
$ws.Range("D106").Value = 33
$ws.Range("E106").Value = 175
$ws.Range("F106").Value = 2468
$ws.Range("G106").Value = 177
$ws.Range("H106").Value = 3000

# --- Row 107 : 17-04-2018 ---
$ws.Range("A107").Value = "17-04-2018"
$ws.Range("B107").Value = 1857
$ws.Range("C107").Value = 76
$ws.Range("D107").Value = 14
$ws.Range("E107").Value = 130
$ws.Range("F107").Value = 3506
$ws.Range("G107").Value = 155
$ws.Range("H107").Value = 3000

# --- Row 108 : 18-04-2018 ---
$ws.Range("A108").Value = "18-04-2018"
$ws.Range("B108").Value = 1812
$ws.Range("C108").Value = 67
$ws.Range("D108").Value = 19
$ws.Range("E108").Value = 165
$ws.Range("F108").Value = 3077
$ws.Range("G108").Value = 127
$ws.Range("H108").Value = 3250

# --- Row 109 : 19-04-2018 ---
$ws.Range("A109").Value = "19-04-2018"
$ws.Range("B109").Value = 1866
$ws.Range("C109").Value = 63
$ws.Range("D109").Value = 27
$ws.Range("E109").Value = 149
$ws.Range("F109").Value = 4621
$ws.Range("G109").Value = 171
$ws.Range("H109").Value = 2500

# --- Row 110 : 20-04-2018 ---
$ws.Range("A110").Value = "20-04-2018"
$ws.Range("B110").Value = 2607
$ws.Range("C110").Value = 103
$ws.Range("D110").Value = 17
$ws.Range("E110").Value = 277
$ws.Range("F110").Value = 4224
$ws.Range("G110").Value = 140
$ws.Range("H110").Value = 2000

# --- Row 111 : 21-04-2018 ---
$ws.Range("A111").Value = "21-04-2018"
$ws.Range("B111").Value = 2171
$ws.Range("C111").Value = 72
$ws.Range("D111").Value = 37
$ws.Range("E111").Value = 204
$ws.Range("F111").Value = 2328
$ws.Range("G111").Value = 159
$ws.Range("H111").Value = 3250

# Extend the "Yes/No" formulas from columns I and J down through row 111,
# same as dragging the fill handle from I105:J105 down to I111:J111.
for ($r = 106; $r -le 111; $r++) {
    $ws.Range("I$r").Formula = "=IF(H$r>=2200,""Yes"",""No"")"
    $ws.Range("J$r").Formula = "=IF(B$r<=1800,""Yes"",""No"")"
}

# Recalculate so cached formula results are written out.
$excel.Calculate()

# Scroll position / active selection left behind by the editor after the edit.
$excel.ActiveWindow.ScrollRow = 82
$ws.Range("L109").Select()
